# Added SPI Flash ROM (U3) to the CPL (component placement list).
# The whole placement table was regenerated by the CAD tool after the new
# part was added, so every designator's Mid X / Mid Y / Rotation shifts,
# a handful of designators get renumbered (J3->J2, a new C18 appears, Q1/Q2
# are inserted before the resistors, and U3 is appended at the end), and
# the used range grows from 52 to 55 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents of the data block, A2:E43 -> Designator, Mid X, Mid Y, Layer, Rotation
$rows = @(
    @("C1", 6.03, 33.34, "Top", 270),
    @("C2", 6.03, 26.99, "Top", 90),
    @("C3", 11.11, 36.51, "Top", 90),
    @("C4", 8.57, 36.51, "Top", 90),
    @("C5", 10.16, 15.88, "Top", 180),
    @("C6", 15.56, 11.11, "Top", 90),
    @("C7", 13.65, 15.56, "Top", 270),
    @("C8", 10.48, 21.59, "Top", 270),
    @("C9", 16.51, 15.56, "Top", 270),
    @("C10", 25.4, 30.8, "Top", 90),
    @("C11", 23.18, 15.56, "Top", 270),
    @("C12", 14.61, 33.02, "Top", 180),
    @("C13", 39.69, 35.88, "Top", 180),
    @("C14", 23.18, 53.66, "Top", 180),
    @("C15", 12.07, 30.8, "Top", 180),
    @("C16", 7.62, 17.78, "Top", 180),
    @("C17", 40.32, 38.1, "Top", 180),
    @("C18", 26.04, 12.38, "Top", 270),
    @("D1", 14.29, 35.24, "Top", 0),
    @("D2", 7.94, 13.97, "Top", 180),
    @("J2", 3.66, 12.07, "Top", 270),
    @("LED1", 2.22, 21.91, "Top", 0),
    @("LED2", 2.22, 20.64, "Top", 0),
    @("LED3", 2.22, 19.37, "Top", 0),
    @("LED4", 2.22, 18.1, "Top", 0),
    @("Q1", 2.54, 30.16, "Top", 90),
    @("Q2", 9.84, 33.97, "Top", 0),
    @("R1", 22.23, 12.38, "Top", 270),
    @("R2", 23.18, 55.88, "Top", 180),
    @("R3", 6.03, 30.16, "Top", 90),
    @("R4", 23.18, 51.44, "Top", 180),
    @("R5", 20.96, 15.56, "Top", 90),
    @("R6", 18.73, 15.56, "Top", 90),
    @("R7", 25.4, 15.56, "Top", 90),
    @("R8", 8.26, 9.84, "Top", 0),
    @("R9", 8.26, 12.07, "Top", 0),
    @("RN1", 6.67, 20, "Top", 270),
    @("SW1", 15.24, 52.07, "Top", 0),
    @("SW2", 31.12, 52.07, "Top", 0),
    @("U1", 19.37, 25.08, "Top", 270),
    @("U2", 12.7, 12.07, "Top", 0),
    @("U3", 30.16, 13.64, "Top", 0)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}

# Used range grows to row 55 (3 extra blank-but-formatted rows at the bottom).
$lastDataRow = $startRow + $rows.Count - 1
$newLastRow = 55
for ($r = $lastDataRow + 1; $r -le $newLastRow; $r++) {
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 5)).Value = ""
}

# Restore the view to the top of the sheet with the whole-sheet row selection
# that Calc records after a full-table refresh.
$ws.Activate()
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$ws.Rows.Select()
